# Updated cryptos list on Mon Aug 21 20:42:24 UTC 2023 with GitHub Actions
# Refresh latest crypto prices / 1h volume deltas from the coinranking feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.183.55"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.677.08"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.77"
$ws.Range("E5").Value = "  -3.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5251"
$ws.Range("E6").Value = "  -4.80%  "
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06299"
$ws.Range("E9").Value = "  -2.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.41"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07559"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "1.683.67"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.452"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5645"
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.95"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008035"
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("D17").Value = "26.254.70"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.829"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "188.12"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("E21").Value = "  -5.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.192"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.20"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1249"
$ws.Range("E25").Value = "  -6.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.592"
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.06"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06191"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.361"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.284"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.498"
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.441"
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6070"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.745"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.099"
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01613"
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("D40").Value = "1.084.36"
$ws.Range("E40").Value = "  -3.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8699"
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").Value = "1.829.24"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000110"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.27"
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9983"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.043"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05238"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.975"
$ws.Range("E51").Value = "  -2.97%  "
